$d = $word.ActiveDocument

# The author corrected the logged date from 2021/11/29 to 2021/11/30
# (commit message: "11/30"). Replace the whole date string in place so
# the visible text matches the new value.
$d.Content.Find.Execute("日期：2021/11/29", $true, $false, $false, $false, $false,
                         $true, 1, $false, "日期：2021/11/30", 2)
